$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.040076637849155
$ws.Range("D2").Value = 1.041324372107741
$ws.Range("E2").Value = 1.038499969606129
$ws.Range("F2").Value = 1.044211778545213
$ws.Range("I2").Value = 1.039955480000692
$ws.Range("J2").Value = 1.045165198421545
$ws.Range("K2").Value = 1.044104209149681
$ws.Range("L2").Value = 1.04128782659133
$ws.Range("M2").Value = 1.046983465405821
$ws.Range("N2").Value = 1.046649452885428

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.041187323421269
$ws.Range("D3").Value = 1.042163992862165
$ws.Range("E3").Value = 1.039450367936969
$ws.Range("F3").Value = 1.046489524004359
$ws.Range("I3").Value = 1.040317008443732
$ws.Range("J3").Value = 1.045920547586684
$ws.Range("K3").Value = 1.044754358160209
$ws.Range("L3").Value = 1.042047880168505
$ws.Range("M3").Value = 1.049068579795246
$ws.Range("N3").Value = 1.04740587473302

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.04190552382056
$ws.Range("D4").Value = 1.042706768362343
$ws.Range("E4").Value = 1.04006521637119
$ws.Range("F4").Value = 1.047957756210611
$ws.Range("I4").Value = 1.040549216918442
$ws.Range("J4").Value = 1.046408272725401
$ws.Range("K4").Value = 1.045173889816304
$ws.Range("L4").Value = 1.042538943838069
$ws.Range("M4").Value = 1.050411852393241
$ws.Range("N4").Value = 1.047894292497378

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.042207341237998
$ws.Range("D5").Value = 1.042934829546342
$ws.Range("E5").Value = 1.04032367088048
$ws.Range("F5").Value = 1.048573689310993
$ws.Range("I5").Value = 1.040646426590012
$ws.Range("J5").Value = 1.046613066289942
$ws.Range("K5").Value = 1.045349985516798
$ws.Range("L5").Value = 1.042745211104521
$ws.Range("M5").Value = 1.050975174913558
$ws.Range("N5").Value = 1.048099376892271

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.042258011100427
$ws.Range("D6").Value = 1.042973114939219
$ws.Range("E6").Value = 1.04036706494821
$ws.Range("F6").Value = 1.04867703124994
$ws.Range("I6").Value = 1.04066272449428
$ws.Range("J6").Value = 1.046647437657733
$ws.Range("K6").Value = 1.045379536663609
$ws.Range("L6").Value = 1.042779834019955
$ws.Range("M6").Value = 1.051069678690479
$ws.Range("N6").Value = 1.048133797071346

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.041909557163746
$ws.Range("D7").Value = 1.042709816202368
$ws.Range("E7").Value = 1.040068669959577
$ws.Range("F7").Value = 1.047965991447241
$ws.Range("I7").Value = 1.040550517449898
$ws.Range("J7").Value = 1.046411010150987
$ws.Range("K7").Value = 1.045176243895057
$ws.Range("L7").Value = 1.042541700679978
$ws.Range("M7").Value = 1.050419384952543
$ws.Range("N7").Value = 1.047897033810421

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.040452101777686
$ws.Range("D8").Value = 1.04160823360937
$ws.Range("E8").Value = 1.038821187599478
$ws.Range("F8").Value = 1.044982737870743
$ws.Range("I8").Value = 1.040078018595855
$ws.Range("J8").Value = 1.045420687992709
$ws.Range("K8").Value = 1.044324171008278
$ws.Range("L8").Value = 1.04154484488222
$ws.Range("M8").Value = 1.047689387878882
$ws.Range("N8").Value = 1.046905305281088

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03788002500426
$ws.Range("D9").Value = 1.03966307897578
$ws.Range("E9").Value = 1.036621948497621
$ws.Range("F9").Value = 1.039681268108246
$ws.Range("I9").Value = 1.039232111189705
$ws.Range("J9").Value = 1.043667585701017
$ws.Range("K9").Value = 1.042813753439985
$ws.Range("L9").Value = 1.039782498733305
$ws.Range("M9").Value = 1.042831883695422
$ws.Range("N9").Value = 1.045149713383058

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.036162542505394
$ws.Range("D10").Value = 1.038363488114897
$ws.Range("E10").Value = 1.035154997047915
$ws.Range("F10").Value = 1.036114796716072
$ws.Range("I10").Value = 1.038659091595444
$ws.Range("J10").Value = 1.042493321584806
$ws.Range("K10").Value = 1.041800663424219
$ws.Range("L10").Value = 1.038603617437449
$ws.Range("M10").Value = 1.039559985185463
$ws.Range("N10").Value = 1.043973781677112

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.035418158047441
$ws.Range("D11").Value = 1.037800054986999
$ws.Range("E11").Value = 1.03451957745557
$ws.Range("F11").Value = 1.034562358698917
$ws.Range("I11").Value = 1.038408782638255
$ws.Range("J11").Value = 1.041983510366109
$ws.Range("K11").Value = 1.0413604982412
$ws.Range("L11").Value = 1.038092177874056
$ws.Range("M11").Value = 1.038134800024821
$ws.Range("N11").Value = 1.043463246468008

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.03514155130806
$ws.Range("D12").Value = 1.03759066297951
$ws.Range("E12").Value = 1.03428351896046
$ws.Range("F12").Value = 1.033984449858918
$ws.Range("I12").Value = 1.038315475297466
$ws.Range("J12").Value = 1.041793938654699
$ws.Range("K12").Value = 1.041196774792648
$ws.Range("L12").Value = 1.037902057433077
$ws.Range("M12").Value = 1.03760411643116
$ws.Range("N12").Value = 1.043273405543024

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.035200889392229
$ws.Range("D13").Value = 1.037635583196279
$ws.Range("E13").Value = 1.034334155947195
$ws.Range("F13").Value = 1.034108471165389
$ws.Range("I13").Value = 1.038335505082618
$ws.Range("J13").Value = 1.041834611757148
$ws.Range("K13").Value = 1.041231904319524
$ws.Range("L13").Value = 1.037942845694582
$ws.Range("M13").Value = 1.037718009622812
$ws.Range("N13").Value = 1.043314136405943

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.035395295896666
$ws.Range("D14").Value = 1.037782748792849
$ws.Range("E14").Value = 1.0345000655218
$ws.Range("F14").Value = 1.034514614612179
$ws.Range("I14").Value = 1.038401076602039
$ws.Range("J14").Value = 1.041967844509817
$ws.Range("K14").Value = 1.04134696944545
$ws.Range("L14").Value = 1.038076465512744
$ws.Range("M14").Value = 1.038090960434683
$ws.Range("N14").Value = 1.043447558364403

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.035515061581761
$ws.Range("D15").Value = 1.037873408008207
$ws.Range("E15").Value = 1.0346022831297
$ws.Range("F15").Value = 1.034764684187993
$ws.Range("I15").Value = 1.038441433375256
$ws.Range("J15").Value = 1.042049906342115
$ws.Range("K15").Value = 1.041417834775916
$ws.Range("L15").Value = 1.038158773275542
$ws.Range("M15").Value = 1.038320573434893
$ws.Range("N15").Value = 1.043529736733915

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.036211929731985
$ws.Range("D16").Value = 1.03840086632527
$ws.Range("E16").Value = 1.03519716292251
$ws.Range("F16").Value = 1.036217651968789
$ws.Range("I16").Value = 1.03867565744883
$ws.Range("J16").Value = 1.042527127461795
$ws.Range("K16").Value = 1.041829844116172
$ws.Range("L16").Value = 1.038637539229849
$ws.Range("M16").Value = 1.03965438902154
$ws.Range("N16").Value = 1.044007635562324

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.036648865947264
$ws.Range("D17").Value = 1.038731537422954
$ws.Range("E17").Value = 1.03557025479105
$ws.Range("F17").Value = 1.037126853271912
$ws.Range("I17").Value = 1.038821992372419
$ws.Range("J17").Value = 1.042826113085447
$ws.Range("K17").Value = 1.042087885837799
$ws.Range("L17").Value = 1.038937593463384
$ws.Range("M17").Value = 1.040488771182284
$ws.Range("N17").Value = 1.044307045779855

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.036903655844724
$ws.Range("D18").Value = 1.038924344743001
$ws.Range("E18").Value = 1.035787851905792
$ws.Range("F18").Value = 1.037656393256456
$ws.Range("I18").Value = 1.038907136167319
$ws.Range("J18").Value = 1.043000376553298
$ws.Range("K18").Value = 1.04223825357136
$ws.Range("L18").Value = 1.039112515996069
$ws.Range("M18").Value = 1.040974641371489
$ws.Range("N18").Value = 1.044481556721819

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.036990521254968
$ws.Range("D19").Value = 1.038990075673603
$ws.Range("E19").Value = 1.035862043395754
$ws.Range("F19").Value = 1.037836821381431
$ws.Range("I19").Value = 1.038936132313426
$ws.Range("J19").Value = 1.043059773959926
$ws.Range("K19").Value = 1.042289500795983
$ws.Range("L19").Value = 1.039172144181234
$ws.Range("M19").Value = 1.041140174212024
$ws.Range("N19").Value = 1.044541038479577

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.036601993844994
$ws.Range("D20").Value = 1.038696066546898
$ws.Range("E20").Value = 1.035530227759684
$ws.Range("F20").Value = 1.037029385742117
$ws.Range("I20").Value = 1.0388063138534
$ws.Range("J20").Value = 1.042794048202028
$ws.Range("K20").Value = 1.042060215278787
$ws.Range("L20").Value = 1.038905410226368
$ws.Range("M20").Value = 1.040399334023512
$ws.Range("N20").Value = 1.044274935360625

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.03533805107009
$ws.Range("D21").Value = 1.037739415189841
$ws.Range("E21").Value = 1.034451210282449
$ws.Range("F21").Value = 1.034395050733584
$ws.Range("I21").Value = 1.038381776596029
$ws.Range("J21").Value = 1.041928616469592
$ws.Range("K21").Value = 1.041313091907189
$ws.Range("L21").Value = 1.038037121939366
$ws.Range("M21").Value = 1.037981172099951
$ws.Range("N21").Value = 1.043408274615861

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.034542726436106
$ws.Range("D22").Value = 1.037137305242152
$ws.Range("E22").Value = 1.033772583363473
$ws.Range("F22").Value = 1.03273140370623
$ws.Range("I22").Value = 1.038112934100528
$ws.Range("J22").Value = 1.041383297005832
$ws.Range("K22").Value = 1.04084203401678
$ws.Range("L22").Value = 1.037490331275426
$ws.Range("M22").Value = 1.036453199345005
$ws.Range("N22").Value = 1.042862180735914

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.03496440433654
$ws.Range("D23").Value = 1.037456555211474
$ws.Range("E23").Value = 1.03413235662877
$ws.Range("F23").Value = 1.033614044207692
$ws.Range("I23").Value = 1.038255635438455
$ws.Range("J23").Value = 1.04167249472397
$ws.Range("K23").Value = 1.04109187592757
$ws.Range("L23").Value = 1.037780277985531
$ws.Range("M23").Value = 1.037263938365725
$ws.Range("N23").Value = 1.043151789147983

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.036623173539706
$ws.Range("D24").Value = 1.038712094517966
$ws.Range("E24").Value = 1.035548314317663
$ws.Range("F24").Value = 1.037073429537284
$ws.Range("I24").Value = 1.038813398952981
$ws.Range("J24").Value = 1.042808537343936
$ws.Range("K24").Value = 1.042072718857397
$ws.Range("L24").Value = 1.038919952736687
$ws.Range("M24").Value = 1.040439749329574
$ws.Range("N24").Value = 1.044289445078776

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.038545442926084
$ws.Range("D25").Value = 1.040166436745428
$ws.Range("E25").Value = 1.037190635983209
$ws.Range("F25").Value = 1.041057331665773
$ws.Range("I25").Value = 1.039452389025458
$ws.Range("J25").Value = 1.044121769436275
$ws.Range("K25").Value = 1.04320530649318
$ws.Range("L25").Value = 1.040238801297047
$ws.Range("M25").Value = 1.044093429798098
$ws.Range("N25").Value = 1.045604542111316
